$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at position 13 (old rows 13-25 shift down to 14-26,
#    carrying their row heights with them).
$ws.Rows.Item(13).Insert()

# Excel copies formatting from the row above into the freshly inserted row;
# remove the stray A13 cell that results (final layout has no A13 cell).
$ws.Range("A13").Clear()

# 2. Populate the new row 13 (B13/C13) with the "Docentes responsaveis" value
#    that used to live in row 13 (now pushed down structurally one level in the
#    shared-string table, but identical text).
$ws.Range("B13").Value = "5840601 - Hélcio José Izário Filho"
$ws.Range("C13").Value = "5840601 - Hélcio José Izário Filho"
$ws.Range("B13").Font.Bold = $false
$ws.Range("B13").WrapText = $true
$ws.Range("B13").VerticalAlignment = -4160
$ws.Range("C13").Font.Bold = $false
$ws.Range("C13").Font.Color = 255
$ws.Range("C13").WrapText = $true
$ws.Range("C13").VerticalAlignment = -4160

# 3. Fill in the new/updated Portuguese + corrected content for the rows that
#    changed text (objectives, short syllabus, full syllabus, method/criteria/
#    recovery shifted down one row, and the new bibliography text).
$ws.Range("B10").Value = "Gerais: - Mostrar a Química Analítica por via úmida como uma ciência que se propõe a determinar a composição qualitativa e quantitativa da matéria por meio de reações químicas específicas e observação crítica dos resultados, requerendo para isso observadores competentes tecnicamente, criativos e sensatos.`nEspecíficos: - Ao concluir o curso os alunos devem: interpretar adequadamente as técnicas e princípios inseridos nos textos de Química Analítica; saber manusear com precisão e eficiência a instrumentação analítica, produtos tóxicos, inflamáveis e cáusticos; compreender os diversos tipos de cálculos estequiométricos; preparar, aferir, conservar e usar adequadamente soluções padrões, bem como, o descarte adequadamente em função da toxicidade dos reagentes/produtos."
$ws.Range("C10").Value = "Gerais: - Mostrar a Química Analítica por via úmida como uma ciência que se propõe a determinar a composição qualitativa e quantitativa da matéria por meio de reações químicas específicas e observação crítica dos resultados, requerendo para isso observadores competentes tecnicamente, criativos e sensatos.`nEspecíficos: - Ao concluir o curso os alunos devem: interpretar adequadamente as técnicas e princípios inseridos nos textos de Química Analítica; saber manusear com precisão e eficiência a instrumentação analítica, produtos tóxicos, inflamáveis e cáusticos; compreender os diversos tipos de cálculos estequiométricos; preparar, aferir, conservar e usar adequadamente soluções padrões, bem como, o descarte adequadamente em função da toxicidade dos reagentes/produtos."
$ws.Range("B14").Value = "- Bases teóricas da química analítica; Introdução à análise qualitativa; Leis e teorias fundamentais; Análise qualitativa sistemática de cátions; Análise qualitativa de ânions.`n- Fundamentos da análise titrimétrica; Titrimetria por Neutralização; Titrimetria por Precipitação; Titrimetria por oxidação-redução: Permanganatometria e Tiossulfatometria; Titrimetria por Complexação."
$ws.Range("C14").Value = "- Bases teóricas da química analítica; Introdução à análise qualitativa; Leis e teorias fundamentais; Análise qualitativa sistemática de cátions; Análise qualitativa de ânions.`n- Fundamentos da análise titrimétrica; Titrimetria por Neutralização; Titrimetria por Precipitação; Titrimetria por oxidação-redução: Permanganatometria e Tiossulfatometria; Titrimetria por Complexação."
$ws.Range("B16").Value = "- Bases Teóricas da Análise Qualitativa: Equilíbrio Químico; Efeito do íon Comum; Produto Iônico da água; Concentração do íon H+; Soluções tampão; Hidrólise dos sais; Produto de solubilidade. Operações analíticas: Precipitação; Filtração, Centrifugação. Equipamentos para ensaios por via úmida. Limpeza da aparelhagem/vidrarias. Classificação analítica dos cátions e dos ânions. Análise Qualitativa Sistemática: separação e identificação dos cátions do 1º, 2º e 3º Grupos. Identificação dos ânions segundo Vogel.`n- Fundamentos de Análise Quantitativa - Titrimetria por Neutralização: fundamentos específicos. Preparação e aferição das soluções padrão ácidas e alcalinas; determinações alcalimétricas e acidimétricas. - Titrimetria por Precipitação: discussão geral da Argentimetria. Princípio da acão dos indicadores. Preparação e emprego da solução padrão de nitrato de prata. Sulfocianetometria. Discussão geral. Preparação, aferição e emprego do processo. - Titrimetria por Oxidação-redução: - Permanganatometria. Características gerais do processo. Preparação / aferição e emprego do processo. - Tiossulfatometria: Características gerais do método. Emprego dos processos titulométricos: direto, inverso, indireto e de retorno. - Complexometria: Estudo teórico da formação de complexos. Grupos de coordenação. Características gerais do método. Preparação, aferição e emprego da solução padrão de EDTAH2Na2. Uso de indicadores metalocrômicos. Determinação de metais bivalentes e trivalentes."
$ws.Range("C16").Value = "- Bases Teóricas da Análise Qualitativa: Equilíbrio Químico; Efeito do íon Comum; Produto Iônico da água; Concentração do íon H+; Soluções tampão; Hidrólise dos sais; Produto de solubilidade. Operações analíticas: Precipitação; Filtração, Centrifugação. Equipamentos para ensaios por via úmida. Limpeza da aparelhagem/vidrarias. Classificação analítica dos cátions e dos ânions. Análise Qualitativa Sistemática: separação e identificação dos cátions do 1º, 2º e 3º Grupos. Identificação dos ânions segundo Vogel.`n- Fundamentos de Análise Quantitativa - Titrimetria por Neutralização: fundamentos específicos. Preparação e aferição das soluções padrão ácidas e alcalinas; determinações alcalimétricas e acidimétricas. - Titrimetria por Precipitação: discussão geral da Argentimetria. Princípio da acão dos indicadores. Preparação e emprego da solução padrão de nitrato de prata. Sulfocianetometria. Discussão geral. Preparação, aferição e emprego do processo. - Titrimetria por Oxidação-redução: - Permanganatometria. Características gerais do processo. Preparação / aferição e emprego do processo. - Tiossulfatometria: Características gerais do método. Emprego dos processos titulométricos: direto, inverso, indireto e de retorno. - Complexometria: Estudo teórico da formação de complexos. Grupos de coordenação. Características gerais do método. Preparação, aferição e emprego da solução padrão de EDTAH2Na2. Uso de indicadores metalocrômicos. Determinação de metais bivalentes e trivalentes."
$ws.Range("B19").Value = "Serão aplicadas, por bimestre, duas avaliações, sendo uma avaliação teórica (peso 0,6) e uma avaliação prática (peso 0,4)."
$ws.Range("C19").Value = "Serão aplicadas, por bimestre, duas avaliações, sendo uma avaliação teórica (peso 0,6) e uma avaliação prática (peso 0,4)."
$ws.Range("B20").Value = "A composição da média P1 e P2 será calculado pelo valor da avaliação teórica x 0,6 mais o valor da avaliação prática x 0,4. A média final será a média aritmética da P1 e P2."
$ws.Range("C20").Value = "A composição da média P1 e P2 será calculado pelo valor da avaliação teórica x 0,6 mais o valor da avaliação prática x 0,4. A média final será a média aritmética da P1 e P2."
$ws.Range("B21").Value = "Na semana da recuperação será dado uma aula teórica e uma avaliação teórica no valor de 10. A Nota final será a média entre a média final (P1 e P2) e a nota da recuperação."
$ws.Range("C21").Value = "Na semana da recuperação será dado uma aula teórica e uma avaliação teórica no valor de 10. A Nota final será a média entre a média final (P1 e P2) e a nota da recuperação."
$ws.Range("B22").Value = "Bibliografia Básica: `n1) VOGEL, Arthur Israel. Química analítica qualitativa. Sao Paulo: Mestrejou, 1981.`n2) VOGEL, Arthur I. Análise química quantitativa/ G. H. Jeffery; J. Bassett; J. Mendham; R. C. Denney. Rio de Janeiro: Guanabara Koogan, 1992.`n3) BACCAN, Nivaldo; ANDRADE, João Carlos de; GODINHO, Oswaldo E.S.; BARONE, José Salvador. Química analítica quantitativa elementar. São Paulo: Edgard Blücher - Instituto Mauá de Tecnologia, 2005-2007.`n4) BACCAN, Nivaldo et al.  Introdução à semimicroanálise qualitativa. Campinas:Editora da UNICAMP, 1988.`nBibliografia Complementar:`nSKOOG, Douglas A. et al. Fundamentos da química analítica. São Paulo: Editora Thomson Learning, 2006-9. `nALEXEYEV, V. Análise Qualitativa. Porto: Editora Lopes da Silva, 1982. `nHARRIS, Daniel C. Análise Química Quantitativa. 6. ed. Rio de Janeiro: Livros Técnicos e Científicos, 2005."
$ws.Range("C22").Value = "Bibliografia Básica: `n1) VOGEL, Arthur Israel. Química analítica qualitativa. Sao Paulo: Mestrejou, 1981.`n2) VOGEL, Arthur I. Análise química quantitativa/ G. H. Jeffery; J. Bassett; J. Mendham; R. C. Denney. Rio de Janeiro: Guanabara Koogan, 1992.`n3) BACCAN, Nivaldo; ANDRADE, João Carlos de; GODINHO, Oswaldo E.S.; BARONE, José Salvador. Química analítica quantitativa elementar. São Paulo: Edgard Blücher - Instituto Mauá de Tecnologia, 2005-2007.`n4) BACCAN, Nivaldo et al.  Introdução à semimicroanálise qualitativa. Campinas:Editora da UNICAMP, 1988.`nBibliografia Complementar:`nSKOOG, Douglas A. et al. Fundamentos da química analítica. São Paulo: Editora Thomson Learning, 2006-9. `nALEXEYEV, V. Análise Qualitativa. Porto: Editora Lopes da Silva, 1982. `nHARRIS, Daniel C. Análise Química Quantitativa. 6. ed. Rio de Janeiro: Livros Técnicos e Científicos, 2005."
